# Appended data to 'Sheet1': a new column E ("CITY") with a header cell
# styled like the existing header row, plus the per-row city values.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header cell E1 ---------------------------------------------------
$ws.Range("E1").Value = "CITY"

# Match the formatting already used by the other header cells (A1:D1) by
# copying their format onto E1 (mirrors Excel's "paste formatting only").
$ws.Range("A1").Copy() | Out-Null
$ws.Range("E1").PasteSpecial(-4122) | Out-Null   # xlPasteFormats

# --- Data rows ----------------------------------------------------------
$ws.Range("E2").Value = "Chennai"

$excel.CutCopyMode = $false
